# Commit: "delete a comment in the ppt file"
#
# The presentation has a single threaded comment (with one reply) attached
# to slide 10 (ppt/comments/comment1.xml, referenced from
# ppt/slides/_rels/slide10.xml.rels). Removing it entirely - including its
# reply - drops the comment part, its relationship, and its
# [Content_Types].xml override, matching the target diff exactly.

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    $commentCount = $s.Comments.Count
    for ($i = $commentCount; $i -ge 1; $i--) {
        $cm = $s.Comments.Item($i)

        # Delete any replies first (a comment with outstanding replies stays
        # behind as an empty thread otherwise).
        $replyCount = $cm.Replies.Count
        for ($j = $replyCount; $j -ge 1; $j--) {
            $cm.Replies.Item($j).Delete()
        }

        $cm.Delete()
    }
}
